$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset previously included an "ECs" sending-cluster block (rows 2-4).
# That block is removed entirely; the FAPs/MuSCs blocks shift up and are
# refreshed with values recomputed from the updated TPM input.
$ws.Range("A2:T4").EntireRow.Delete()

$newRows = @(
    @("FAPs", "Fgf18", "Fgfr3", "ECs", 3, 1, 8.770835333333332, 26.312506, 0.9145494540267, 0.9145494540267, 3, 1, 4.959409333333333, 14.878228, 0.8271666313262851, 0.8271666313262852, 43.4981626132631, 391.483463519368, 0.7564847910685586, 0.7564847910685587),
    @("FAPs", "Fgf18", "Fgfr3", "FAPs", 3, 1, 8.770835333333332, 26.312506, 0.9145494540267, 0.9145494540267, 2, 0.6666666666666666, 0.5648773333333333, 1.694632, 0.09421438109281059, 0.09421438109281059, 4.954446074199111, 44.590014667792, 0.08616371078989338, 0.08616371078989338),
    @("FAPs", "Fgf18", "Fgfr3", "MuSCs", 3, 1, 8.770835333333332, 26.312506, 0.9145494540267, 0.9145494540267, 3, 1, 0.4713726666666667, 1.414118, 0.07861898758090437, 0.07861898758090438, 4.134332039967555, 37.208988359708, 0.07190095216824799, 0.071900952168248),
    @("MuSCs", "Fgf18", "Fgfr3", "ECs", 3, 1, 0.8194993333333332, 2.458498, 0.08545054597330007, 0.08545054597330005, 3, 1, 4.959409333333333, 14.878228, 0.8271666313262851, 0.8271666313262852, 4.064232642393777, 36.57809378154399, 0.07068184025772647, 0.07068184025772647),
    @("MuSCs", "Fgf18", "Fgfr3", "FAPs", 3, 1, 0.8194993333333332, 2.458498, 0.08545054597330007, 0.08545054597330005, 2, 0.6666666666666666, 0.5648773333333333, 1.694632, 0.09421438109281059, 0.09421438109281059, 0.4629165980817778, 4.166249382735999, 0.008050670302917224, 0.008050670302917224),
    @("MuSCs", "Fgf18", "Fgfr3", "MuSCs", 3, 1, 0.8194993333333332, 2.458498, 0.08545054597330007, 0.08545054597330005, 3, 1, 0.4713726666666667, 1.414118, 0.07861898758090437, 0.07861898758090438, 0.3862895860848888, 3.476606274763999, 0.006718035412656376, 0.006718035412656376)
)

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowData = $newRows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}
